$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the entire contents of data rows 2 and 3 (the two
# occurrence records change places), including the stray empty "L" cell
# that only one of the two rows carries.
#
# NOTE: in this interop, Range.Value's *getter* is unreliable (it can
# return a stray reflection/member-description object instead of the
# actual cell content), so values are read back via Value2, which works
# correctly for both get and set.

# --- Capture current (row 2 / row 3) values for the columns that differ ---
$row2 = @{
    A = $ws.Range("A2").Value2
    B = $ws.Range("B2").Value2
    E = $ws.Range("E2").Value2
    F = $ws.Range("F2").Value2
    G = $ws.Range("G2").Value2
    H = $ws.Range("H2").Value2
    I = $ws.Range("I2").Value2
    Q = $ws.Range("Q2").Value2
    R = $ws.Range("R2").Value2
    S = $ws.Range("S2").Value2
}

$row3 = @{
    A = $ws.Range("A3").Value2
    B = $ws.Range("B3").Value2
    E = $ws.Range("E3").Value2
    F = $ws.Range("F3").Value2
    G = $ws.Range("G3").Value2
    H = $ws.Range("H3").Value2
    I = $ws.Range("I3").Value2
    Q = $ws.Range("Q3").Value2
    R = $ws.Range("R3").Value2
    S = $ws.Range("S3").Value2
}

# Column "I" ("Antal") is stored as text in the sheet (e.g. "1", "20"),
# even though it looks numeric. Force it to stay text on write, otherwise
# a plain Value2 assignment of a numeric-looking string is auto-converted
# to a real number. Apply the text format to both destination cells at
# once so they share a single style entry.
$ws.Range("I2:I3").NumberFormat = "@"

# --- Write row 3's old values into row 2 ---
$ws.Range("A2").Value2 = $row3.A
$ws.Range("B2").Value2 = $row3.B
$ws.Range("E2").Value2 = $row3.E
$ws.Range("F2").Value2 = $row3.F
$ws.Range("G2").Value2 = $row3.G
$ws.Range("H2").Value2 = $row3.H
$ws.Range("I2").Value2 = [string]$row3.I
$ws.Range("Q2").Value2 = $row3.Q
$ws.Range("R2").Value2 = $row3.R
$ws.Range("S2").Value2 = $row3.S

# --- Write row 2's old values into row 3 ---
$ws.Range("A3").Value2 = $row2.A
$ws.Range("B3").Value2 = $row2.B
$ws.Range("E3").Value2 = $row2.E
$ws.Range("F3").Value2 = $row2.F
$ws.Range("G3").Value2 = $row2.G
$ws.Range("H3").Value2 = $row2.H
$ws.Range("I3").Value2 = [string]$row2.I
$ws.Range("Q3").Value2 = $row2.Q
$ws.Range("R3").Value2 = $row2.R
$ws.Range("S3").Value2 = $row2.S

# Restore the default ("Normal") style on the I column cells so only their
# content type changes to text - not their visible formatting/style index.
$ws.Range("I2:I3").Style = "Normal"

# --- The empty "L" placeholder cell moves from row 2 to row 3 ---
$ws.Range("L2").ClearContents()
# Touch L3 (which currently has no recorded cell at all) with a no-op
# formatting operation so it becomes a real, but still value-less, blank
# cell - mirroring the blank placeholder cell that used to live at L2.
$ws.Range("L3").ClearFormats()
